# Apply Vietnamese-text revisions to the QA document.
#
# Each substitution is performed with Word's Find & Replace engine
# (wdReplaceAll) scoped to the specific paragraph's Range, so that:
#   - identical substrings elsewhere in the document are left untouched
#   - MatchWholeWord is avoided for strings beginning/ending on spaces or
#     punctuation (whole-word matching requires word-character boundaries)
#   - every Find What string sits fully inside plain (non spell-checked)
#     runs, so no <w:proofErr> spellStart/spellEnd pairs get split

$d = $word.ActiveDocument

function Replace-In-Paragraph($paraIndex, $findText, $replaceText) {
    $rng = $d.Paragraphs($paraIndex).Range
    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, `
                       $true, 1, $false, $replaceText, 2) | Out-Null
}

# Paragraph 6: "Cách sử dụng tools like Jira and Zephyr to write test cases and defects reporting."
#           -> "Cách sử dụng phần mềm Jira and Zephyr để viết test case và test report."
Replace-In-Paragraph 6 "tools like" "phần mềm"
Replace-In-Paragraph 6 "to write test cases and defects reporting." "để viết test case và test report."

# Paragraph 7: "Manual Test một số package release của Iritech (IriTracker win32, ... IriCoreLicense…)"
#           -> "Manual Test một số phần mềm đã hoàn thành của Iritech như: IriTracker win32, ... IriCoreLicense…"
Replace-In-Paragraph 7 "package release" "phần mềm đã hoàn thành"
Replace-In-Paragraph 7 " (" " như: "
Replace-In-Paragraph 7 "…)" "…"

# Paragraph 12: "Sử dụng công cụ lập trình Visual Studio 2008, 2019, 2022, Android Studio..."
#            -> "Sử dụng công cụ lập trình Visual Studio, Android Studio..."
Replace-In-Paragraph 12 "Visual Studio 2008, 2019, 2022, Android Studio" "Visual Studio, Android Studio"
